$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "wake"
$ws.Range("C3").Value = "stage1"
$ws.Range("C4").Value = "stage2"
$ws.Range("C5").Value = "sws"
$ws.Range("C6").Value = "sws"
$ws.Range("C7").Value = "rem"

$ws.Range("C8").Select()
